$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("133:133").Insert()

$ws.Range("A133").Value = 10
$ws.Range("B133").Value = "Vega Modelo de Temuco"
$ws.Range("C133").Value = "La Araucanía"
$ws.Range("D133").Value = 44781
$ws.Range("E133").Value = 9
$ws.Range("F133").Value = 100112039
$ws.Range("G133").Value = "Ciboulette"
$ws.Range("H133").Value = "Sin especificar"
$ws.Range("I133").Value = "Primera"
$ws.Range("J133").Value = 90
$ws.Range("K133").Value = 4000
$ws.Range("L133").Value = 5000
$ws.Range("M133").Value = 4556
$ws.Range("N133").Value = "`$/docena de atados"
$ws.Range("O133").Value = "Provincia de Cautín"
$ws.Range("P133").Value = 1519
$ws.Range("Q133").Value = 3
$ws.Range("R133").Value = "Hortaliza"
